$ws = $excel.ActiveWorkbook.ActiveSheet
$win = $excel.ActiveWindow
Write-Host "initial ScrollColumn:" $win.ScrollColumn  "ScrollRow:" $win.ScrollRow
$ws.Range("E17").Select()
Write-Host "after select ScrollColumn:" $win.ScrollColumn  "ScrollRow:" $win.ScrollRow
$win.ScrollColumn = 3
Write-Host "after set ScrollColumn:" $win.ScrollColumn  "ScrollRow:" $win.ScrollRow
